$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in column B for rows 1-6 with tracklet labels (version 2 of tracklet data)
$ws.Range("B1").Value = 0
$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 1
$ws.Range("B4").Value = 1
$ws.Range("B5").Value = 1
$ws.Range("B6").Value = 1

# Update the view: scroll position and active selection
$ws.Range("B1").Select()
$excel.ActiveWindow.ScrollRow = 25
